$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ 'C'=9.628476381622647; 'D'=4.889477451973889; 'E'=12.53729045262738; 'F'=30.06386191256705; 'G'=3.626191705832343; 'I'=26.77851708482352; 'L'=9.461361224737708; 'O'=26.27224278197679 }
    3 = @{ 'C'=9.616662997744301; 'D'=4.904063324356488; 'E'=12.50587338210542; 'F'=29.61645569846391; 'G'=3.629809073164946; 'I'=26.52703300356137; 'L'=9.462156847846751; 'O'=25.99078460847631 }
    4 = @{ 'C'=9.611393688489276; 'D'=4.913540692465935; 'E'=12.48947518498029; 'F'=29.34776987289026; 'G'=3.632147283375773; 'I'=26.37914352774891; 'L'=9.464457197895868; 'O'=25.82415151807154 }
    5 = @{ 'C'=9.609746457278014; 'D'=4.917534230268341; 'E'=12.48352350784806; 'F'=29.23993525346502; 'G'=3.633129687416996; 'I'=26.32057785934655; 'L'=9.465850124014379; 'O'=25.7578749704657 }
    6 = @{ 'C'=9.609503163802815; 'D'=4.918205300668106; 'E'=12.48257945799152; 'F'=29.22213357267391; 'G'=3.633294603725373; 'I'=26.31095746928115; 'L'=9.466108923981693; 'O'=25.74697017307598 }
    7 = @{ 'C'=9.611369447427201; 'D'=4.913594018198493; 'E'=12.4893919558401; 'F'=29.34630868273633; 'G'=3.632160412569632; 'I'=26.37834672866234; 'L'=9.46447413928597; 'O'=25.82325100505348 }
    8 = @{ 'C'=9.623992090321806; 'D'=4.894398571535361; 'E'=12.52585992806422; 'F'=29.9084324647919; 'G'=3.627414729997184; 'I'=26.69048884695129; 'L'=9.461259512296792; 'O'=26.17395749425673 }
    9 = @{ 'C'=9.664425100048851; 'D'=4.860882506918765; 'E'=12.62013471018179; 'F'=31.05216523296169; 'G'=3.61903281526939; 'I'=27.3515373107982; 'L'=9.469328790192879; 'O'=26.90738542557767 }
    10 = @{ 'C'=9.703584947228403; 'D'=4.838756717498568; 'E'=12.70298026706907; 'F'=31.90909069558753; 'G'=3.613431099136407; 'I'=27.86307116331799; 'L'=9.484005659677566; 'O'=27.46939495718129 }
    11 = @{ 'C'=9.723422804829166; 'D'=4.829230182980425; 'E'=12.74354063641782; 'F'=32.30072432240082; 'G'=3.611002066040538; 'I'=28.10050732186045; 'L'=9.492575029731206; 'O'=27.7290626869189 }
    12 = @{ 'C'=9.73122278003323; 'D'=4.825699937287994; 'E'=12.75930534164623; 'F'=32.4491386887228; 'G'=3.610099282081606; 'I'=28.19102285860892; 'L'=9.496091200808451; 'O'=27.82788060060402 }
    13 = @{ 'C'=9.729530170363496; 'D'=4.826456806450007; 'E'=12.7558922328186; 'F'=32.41717239182121; 'G'=3.610292956648078; 'I'=28.1715031122685; 'L'=9.495321890515534; 'O'=27.80657812986101 }
    14 = @{ 'C'=9.724058765028294; 'D'=4.828938200548288; 'E'=12.74482953108409; 'F'=32.31293307529683; 'G'=3.610927452589678; 'I'=28.1079424381975; 'L'=9.49285888001239; 'O'=27.7371832362148 }
    15 = @{ 'C'=9.720744751070232; 'D'=4.830468179931331; 'E'=12.73810585488765; 'F'=32.24909359555821; 'G'=3.611318315440208; 'I'=28.06908598802575; 'L'=9.49138549025899; 'O'=27.69473769911101 }
    16 = @{ 'C'=9.702328954494428; 'D'=4.839390119283301; 'E'=12.70038671921716; 'F'=31.88352013745194; 'G'=3.613592232106707; 'I'=27.8476428558882; 'L'=9.483483634649936; 'O'=27.45249809503832 }
    17 = @{ 'C'=9.691547637645554; 'D'=4.845001240094932; 'E'=12.67797786151858; 'F'=31.65960957850508; 'G'=3.615017665516661; 'I'=27.71295294107355; 'L'=9.479120216899565; 'O'=27.3048549978935 }
    18 = @{ 'C'=9.68553721851792; 'D'=4.848279315373789; 'E'=12.6653598592197; 'F'=31.53100161992122; 'G'=3.615848764499557; 'I'=27.63593341792979; 'L'=9.476788728982752; 'O'=27.22031705978993 }
    19 = @{ 'C'=9.683535040638251; 'D'=4.849397930656454; 'E'=12.66113438942697; 'F'=31.48749262216883; 'G'=3.616132091853076; 'I'=27.60993566519988; 'L'=9.476029967427637; 'O'=27.19176236151778 }
    20 = @{ 'C'=9.692675611766747; 'D'=4.844398680275354; 'E'=12.68033533196423; 'F'=31.6834277618024; 'G'=3.614864764470275; 'I'=27.72724484181588; 'L'=9.479566270463376; 'O'=27.32053291759321 }
    21 = @{ 'C'=9.725658065717575; 'D'=4.828207260086159; 'E'=12.74806797866398; 'F'=32.34354884325683; 'G'=3.610740624173945; 'I'=28.12659598749735; 'L'=9.493574976728548; 'O'=27.75755367790955 }
    22 = @{ 'C'=9.748889841475499; 'D'=4.818075393189774; 'E'=12.79469392360504; 'F'=32.77555802606552; 'G'=3.608144523164053; 'I'=28.39107715453985; 'L'=9.504310230473529; 'O'=28.04597191127868 }
    23 = @{ 'C'=9.736338418892061; 'D'=4.823441836120138; 'E'=12.76959570744074; 'F'=32.54498168427403; 'G'=3.609521063088684; 'I'=28.24962569373626; 'L'=9.49843647682386; 'O'=27.89181073488928 }
    24 = @{ 'C'=9.692165069143421; 'D'=4.844670935041286; 'E'=12.67926869310605; 'F'=31.67265917207487; 'G'=3.614933854897314; 'I'=27.72078217329668; 'L'=9.479364057935955; 'O'=27.31344385160795 }
    25 = @{ 'C'=9.651817675125459; 'D'=4.86950954277092; 'E'=12.59222161091889; 'F'=30.7392129463442; 'G'=3.621202115862117; 'I'=27.16786770945215; 'L'=9.465607692995173; 'O'=26.70455883351526 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}

Write-Output "Applied $($data.Count) rows of updates"
